$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the B:H values between row 2 and row 3 (A, I, J stay as they are per row)
$cols = @("B","C","D","E","F","G","H")

foreach ($col in $cols) {
    $v2 = $ws.Range($col + "2").Value2
    $v3 = $ws.Range($col + "3").Value2
    $ws.Range($col + "2").Value = $v3
    $ws.Range($col + "3").Value = $v2
}
